$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, pushing the existing row 17 (and the blank rows /
# signature block below it) down by one.
$ws.Rows.Item(17).Insert()

# The newly inserted row 17 should look like row 16 (same formatting as the
# first data row) before we fill in its values.
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# --- Data updates ---

# VALOR MORA total
$ws.Range("E11").Value = 170820

# Cant. Periodos
$ws.Range("F13").Value = 3

# Row 16 (first data row): Periodo Mora = 2506
$ws.Range("E16").Value = "2506"

# Row 17 (new second data row, period 2507)
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1044922077"
$ws.Range("D17").Value = "SUHAIL BAEZ AYALA"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18 (previously row 17, now shifted down; period 2508)
$ws.Range("E18").Value = "2508"

Write-Output "done"
